$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 430.25
$ws.Range("J12").Value = 612.4
$ws.Range("L12").Value = 612.4
$ws.Range("N12").Value = -952.4
$ws.Range("H21").Value = 8284.75
$ws.Range("J21").Value = 5319.5
$ws.Range("L21").Value = 5319.5
$ws.Range("N21").Value = -6255.5
$ws.Range("H23").Value = 8284.75
$ws.Range("J23").Value = 5319.5
$ws.Range("L23").Value = 5319.5
$ws.Range("N23").Value = -5787.5
$ws.Range("H29").Value = 2652.5
$ws.Range("I29").Value = 42
$ws.Range("J29").Value = 5263
$ws.Range("K29").Value = 126
$ws.Range("L29").Value = 15789
$ws.Range("M29").Value = 155
$ws.Range("N29").Value = -16351
$ws.Range("H38").Value = 1792.0625
$ws.Range("J38").Value = 3627
$ws.Range("L38").Value = 10881
$ws.Range("N38").Value = -11625
$ws.Range("H40").Value = 3879.7666
$ws.Range("I40").Value = 3208.4211
$ws.Range("K40").Value = 3208.4211
$ws.Range("M40").Value = -3033.4211
$ws.Range("H43").Value = 1999.75
$ws.Range("I43").Value = 3000
$ws.Range("J43").Value = 1666.3334
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 1666.3334
$ws.Range("M43").Value = -2931
$ws.Range("N43").Value = -1804.3334
$ws.Range("H58").Value = 1423.3334
$ws.Range("I58").Value = 97
$ws.Range("J58").Value = 2749.6667
$ws.Range("K58").Value = 291
$ws.Range("L58").Value = 8249.000100000001
$ws.Range("M58").Value = -141
$ws.Range("N58").Value = -8549.000100000001
$ws.Range("H70").Value = 8162.2
$ws.Range("I70").Value = 1300
$ws.Range("J70").Value = 8924.666999999999
$ws.Range("K70").Value = 3900
$ws.Range("L70").Value = 26774.001
$ws.Range("M70").Value = -3630
$ws.Range("N70").Value = -27314.001
$ws.Range("H73").Value = 8162.2
$ws.Range("I73").Value = 1300
$ws.Range("J73").Value = 8924.666999999999
$ws.Range("K73").Value = 3900
$ws.Range("L73").Value = 26774.001
$ws.Range("M73").Value = -2964
$ws.Range("N73").Value = -28646.001
$ws.Range("H92").Value = 171.22223
$ws.Range("I92").Value = 162.28572
$ws.Range("K92").Value = 162.28572
$ws.Range("M92").Value = 1085.71428
$ws.Range("H98").Value = 1207.375
$ws.Range("I98").Value = 1049.4
$ws.Range("K98").Value = 1049.4
$ws.Range("M98").Value = 448.5999999999999
$ws.Range("H122").Value = 1207.375
$ws.Range("I122").Value = 1049.4
$ws.Range("K122").Value = 3148.2
$ws.Range("M122").Value = -698.2000000000003

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 19690000
$ws.Range("I6").Value = 21316666
$ws.Range("J6").Value = 17250000
$ws.Range("K6").Value = 21316666
$ws.Range("L6").Value = 17250000
$ws.Range("M6").Value = -21316493
$ws.Range("N6").Value = -17250346
$ws.Range("H61").Value = 2782.4
$ws.Range("I61").Value = 2782.4
$ws.Range("K61").Value = 2782.4
$ws.Range("M61").Value = -2570.4
$ws.Range("H74").Value = 6310.12
$ws.Range("I74").Value = 6539.6113
$ws.Range("J74").Value = 5720
$ws.Range("K74").Value = 6539.6113
$ws.Range("L74").Value = 5720
$ws.Range("M74").Value = -5665.6113
$ws.Range("N74").Value = -7468
$ws.Range("H77").Value = 6310.12
$ws.Range("I77").Value = 6539.6113
$ws.Range("J77").Value = 5720
$ws.Range("K77").Value = 32698.0565
$ws.Range("L77").Value = 28600
$ws.Range("M77").Value = -28330.0565
$ws.Range("N77").Value = -37336
$ws.Range("H97").Value = 961.3158
$ws.Range("I97").Value = 736.8182
$ws.Range("J97").Value = 1270
$ws.Range("K97").Value = 736.8182
$ws.Range("L97").Value = 1270
$ws.Range("M97").Value = -240.8182
$ws.Range("N97").Value = -2262
$ws.Range("H105").Value = 10370
$ws.Range("J105").Value = 10370
$ws.Range("L105").Value = 10370
$ws.Range("N105").Value = -17358
$ws.Range("H132").Value = 2900
$ws.Range("I132").Value = 2900
$ws.Range("K132").Value = 8700
$ws.Range("M132").Value = -6170
$ws.Range("H136").Value = 2782.4
$ws.Range("I136").Value = 2782.4
$ws.Range("K136").Value = 8347.200000000001
$ws.Range("M136").Value = -5797.200000000001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 14000
$ws.Range("I54").Value = 14000
$ws.Range("K54").Value = 14000
$ws.Range("M54").Value = -13516
$ws.Range("H86").Value = 4036.0833
$ws.Range("I86").Value = 1803.25
$ws.Range("J86").Value = 8501.75
$ws.Range("K86").Value = 1803.25
$ws.Range("L86").Value = 8501.75
$ws.Range("M86").Value = -680.25
$ws.Range("N86").Value = -10747.75
$ws.Range("H89").Value = 4036.0833
$ws.Range("I89").Value = 1803.25
$ws.Range("J89").Value = 8501.75
$ws.Range("K89").Value = 9016.25
$ws.Range("L89").Value = 42508.75
$ws.Range("M89").Value = -3400.25
$ws.Range("N89").Value = -53740.75
$ws.Range("H134").Value = 703.6667
$ws.Range("I134").Value = 703.6667
$ws.Range("K134").Value = 2111.0001
$ws.Range("M134").Value = 423.9998999999998

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 87
$ws.Range("I7").Value = 64.5
$ws.Range("K7").Value = 64.5
$ws.Range("M7").Value = 48.5
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").Value = ""
$ws.Range("H134").Value = 2032.2727
$ws.Range("I134").Value = 2065.5
$ws.Range("J134").Value = 1700
$ws.Range("K134").Value = 6196.5
$ws.Range("L134").Value = 5100
$ws.Range("M134").Value = -3661.5
$ws.Range("N134").Value = -10170

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1772.8889
$ws.Range("I34").Value = 674.9091
$ws.Range("K34").Value = 2024.7273
$ws.Range("M34").Value = -1940.7273
$ws.Range("H39").Value = 7501
$ws.Range("J39").Value = 8841.200000000001
$ws.Range("L39").Value = 26523.6
$ws.Range("N39").Value = -27111.6
$ws.Range("H55").Value = 4291.1
$ws.Range("I55").Value = 1529.3334
$ws.Range("J55").Value = 5474.7144
$ws.Range("K55").Value = 4588.0002
$ws.Range("L55").Value = 16424.1432
$ws.Range("M55").Value = -4411.0002
$ws.Range("N55").Value = -16778.1432
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("H131").Value = 3155.5833
$ws.Range("I131").Value = 3260
$ws.Range("J131").Value = 3120.7778
$ws.Range("K131").Value = 9780
$ws.Range("L131").Value = 9362.3334
$ws.Range("M131").Value = -4740
$ws.Range("N131").Value = -19442.3334
$ws.Range("H138").Value = 8833
$ws.Range("I138").Value = 5015
$ws.Range("J138").Value = 9787.5
$ws.Range("K138").Value = 15045
$ws.Range("L138").Value = 29362.5
$ws.Range("M138").Value = -9905
$ws.Range("N138").Value = -39642.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1802.875
$ws.Range("I102").Value = 1631.8572
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 1631.8572
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -9.857199999999921
$ws.Range("N102").Value = -6244
$ws.Range("H122").Value = 3699
$ws.Range("I122").Value = 3748.75
$ws.Range("K122").Value = 11246.25
$ws.Range("M122").Value = -8796.25
$ws.Range("H126").Value = 2993.8
$ws.Range("J126").Value = 8000
$ws.Range("L126").Value = 24000
$ws.Range("N126").Value = -28940
$ws.Range("H132").Value = 95903.63
$ws.Range("I132").Value = 203508.8
$ws.Range("J132").Value = 6232.6665
$ws.Range("K132").Value = 610526.3999999999
$ws.Range("L132").Value = 18697.9995
$ws.Range("M132").Value = -607996.3999999999
$ws.Range("N132").Value = -23757.9995

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 993
$ws.Range("I16").Value = 993
$ws.Range("K16").Value = 993
$ws.Range("M16").Value = -823
$ws.Range("H40").Value = 4200
$ws.Range("I40").Value = 3400
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 3400
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -3264
$ws.Range("N40").Value = -5272
$ws.Range("H46").Value = 4177.353
$ws.Range("J46").Value = 4066.1667
$ws.Range("L46").Value = 4066.1667
$ws.Range("N46").Value = -4442.1667
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = ""
$ws.Range("H136").Value = 2102.5
$ws.Range("I136").Value = 2241.7144
$ws.Range("J136").Value = 1777.6666
$ws.Range("K136").Value = 6725.1432
$ws.Range("L136").Value = 5332.9998
$ws.Range("M136").Value = -4175.1432
$ws.Range("N136").Value = -10432.9998

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 924.75
$ws.Range("I81").Value = 924.75
$ws.Range("K81").Value = 1849.5
$ws.Range("M81").Value = -788.5
$ws.Range("H84").Value = 924.75
$ws.Range("I84").Value = 924.75
$ws.Range("K84").Value = 9247.5
$ws.Range("M84").Value = -3943.5
$ws.Range("H132").Value = 1496.7059
$ws.Range("I132").Value = 1067.5714
$ws.Range("J132").Value = 3499.3333
$ws.Range("K132").Value = 3202.7142
$ws.Range("L132").Value = 10497.9999
$ws.Range("M132").Value = -672.7142000000003
$ws.Range("N132").Value = -15557.9999
$ws.Range("H136").Value = 2723.75
$ws.Range("J136").Value = 7995
$ws.Range("L136").Value = 23985
$ws.Range("N136").Value = -29085
